$d = $word.ActiveDocument
$last = $d.Paragraphs.Last
$r = $last.Range
$insertionPoint = $d.Range($r.Start, $r.Start)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>DROP TABLE PAYMENT;</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>CREATE TABLE PAYMENT</w:t></w:r></w:p><w:p><w:r><w:t>(</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">MARCHENT_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">MARCHENT_BRANCH_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">USER_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">HISTORY_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">CONSTRAINT PK_PAYMENT_ PRIMARY </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>USER_ID,HISTORY_ID,MARCHENT_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_USER_ID_PAYMENT FOREIGN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>USER_ID) REFERENCES USERS(USER_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_HISTORY_ID_PAYMENT FOREIGN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>HISTORY_ID) REFERENCES HISTORY(HISTORY_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_MARCHENT_ID_PAYMENT FOREIGN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>MARCHENT_ID,MARCHENT_BRANCH_ID) REFERENCES MARCHENTS(MARCHENT_ID,MARCHENT_BRANCH_ID)</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>CREATE TABLE CASH_IN</w:t></w:r></w:p><w:p><w:r><w:t>(</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CUSTOMER_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">AGENT_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">HISTORY_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>CONSTRAINT PK_CASH_IN PRIMARY KEY (CUSTOMER_ID</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,AGENT</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_ID,HISTORY_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_CUSTOMER_ID_C_I FOREIGN KEY (CUSTOMER_ID) REFERENCES </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>CUSTOMER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>CUSTOMER_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_AGENT_ID_C_I FOREIGN KEY (AGENT_ID) REFERENCES </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>AGENT(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>AGENT_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_HISTORY_ID_C_I FOREIGN KEY (HISTORY_ID) REFERENCES </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>HISTORY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>HISTORY_ID)</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>CREATE TABLE MOBILE_OPERATOR</w:t></w:r></w:p><w:p><w:r><w:t>(</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">OPERATOR_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">OPERATOR_NAME </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>VARCHAR2(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>255),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">OPERATOR_BANK_AC_NO </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">CONSTRAINT PK_MOBILE_OPERATOR PRIMARY </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>OPERATOR_ID)</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>DROP TABLE MOBILE_RECHARGE;</w:t></w:r></w:p><w:p><w:r><w:t>CREATE TABLE MOBILE_RECHARGE</w:t></w:r></w:p><w:p><w:r><w:t>(</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">OPERATOR_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">USER_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">HISTORY_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">CONSTRAINT PK_MOBILE_RECHARGE PRIMARY </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>OPERATOR_ID,USER_ID,HISTORY_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_OPERATOR_ID_RECHARGE FOREIGN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>OPERATOR_ID) REFERENCES MOBILE_OPERATOR(OPERATOR_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_USER_ID_RECHARGE FOREIGN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>USER_ID) REFERENCES USERS(USER_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_HISTORY_ID_RECHARGE FOREIGN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>HISTORY_ID) REFERENCES HISTORY(HISTORY_ID)</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>CREATE TABLE UTILITY_SERVICE</w:t></w:r></w:p><w:p><w:r><w:t>(</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">SERVICE_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">SERVICE_NAME </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>VARCHAR2(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>255),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">SERVICE_BANK_AC_NO </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">CONSTRAINT PK_UTILITY_SERVICE PRIMARY </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>SERVICE_ID)</w:t></w:r></w:p><w:p><w:r><w:t>);</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>CREATE TABLE PAY_UTILITY_BILL</w:t></w:r></w:p><w:p><w:r><w:t>(</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">SERVICE_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">USER_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">HISTORY_ID </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NUMBER(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>38),</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">CONSTRAINT PK_PAY_UTILITY_BILL PRIMARY </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>SERVICE_ID,USER_ID,HISTORY_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_USER_ID_PAY_U_B FOREIGN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>USER_ID) REFERENCES USERS(USER_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_HISTORY_ID_PAY_U_B FOREIGN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>HISTORY_ID) REFERENCES HISTORY(HISTORY_ID),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CONSTRAINT FK_SERVICE_ID_PAY_U_B FOREIGN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>KEY(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>SERVICE_ID) REFERENCES</w:t></w:r></w:p><w:p><w:r><w:t>UTILITY_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>SERVICE(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>SERVICE_ID)</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xml)
Write-Output "Inserted new SQL schema content."
